# Updates the cryptocurrency price (D) and 1h volume change (E) columns
# for rows 2-51 on the active sheet, matching the latest scrape snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = "25.882.79"; E = "  -0.62%  " },
    @{ Row = 3; D = "1.637.66"; E = "  -0.21%  " },
    @{ Row = 4; D = "1.005"; E = "  -1.44%  " },
    @{ Row = 5; D = "214.87"; E = "  -0.47%  " },
    @{ Row = 6; D = "0.5030"; E = "  +0.47%  " },
    @{ Row = 7; D = "1.004"; E = "  -1.15%  " },
    @{ Row = 8; D = "0.2569"; E = "  -0.32%  " },
    @{ Row = 9; D = "0.06378"; E = "  -0.79%  " },
    @{ Row = 10; D = "19.45"; E = "  -0.14%  " },
    @{ Row = 11; D = "0.07793"; E = "  +0.22%  " },
    @{ Row = 12; D = "4.259"; E = "  -0.04%  " },
    @{ Row = 13; D = "1.644.26"; E = "  +0.19%  " },
    @{ Row = 14; D = "1.862.09"; E = "  -0.14%  " },
    @{ Row = 15; D = "0.5417"; E = "  -0.72%  " },
    @{ Row = 16; D = "0.0₅7887"; E = "  -0.47%  " },
    @{ Row = 17; D = "64.50"; E = "  +1.30%  " },
    @{ Row = 18; D = "25.912.34"; E = "  -0.36%  " },
    @{ Row = 19; D = "1.004"; E = "  -1.36%  " },
    @{ Row = 20; D = "196.76"; E = "  -3.55%  " },
    @{ Row = 21; D = "4.372"; E = "  +1.24%  " },
    @{ Row = 22; D = "9.916"; E = "  -0.98%  " },
    @{ Row = 23; D = "5.964"; E = "  -0.17%  " },
    @{ Row = 24; D = "1.005"; E = "  -1.17%  " },
    @{ Row = 25; D = "1.885"; E = "  -4.68%  " },
    @{ Row = 26; D = "139.89"; E = "  -1.22%  " },
    @{ Row = 27; D = "0.1137"; E = "  -1.44%  " },
    @{ Row = 28; D = "6.828"; E = "  +0.35%  " },
    @{ Row = 29; D = "15.67"; E = "  -0.64%  " },
    @{ Row = 30; D = "1.239"; E = "  -0.28%  " },
    @{ Row = 31; D = "0.04861"; E = "  -3.72%  " },
    @{ Row = 32; D = "3.248"; E = "  -0.53%  " },
    @{ Row = 33; D = "3.176"; E = "  -0.90%  " },
    @{ Row = 34; D = "1.531"; E = "  -1.10%  " },
    @{ Row = 35; D = "2.364"; E = "  +0.01%  " },
    @{ Row = 36; D = "0.8877"; E = "  -0.57%  " },
    @{ Row = 37; D = "2.607"; E = "  -0.51%  " },
    @{ Row = 38; D = "0.5518"; E = "  -2.30%  " },
    @{ Row = 39; D = "1.126.45"; E = "  -0.11%  " },
    @{ Row = 40; D = "0.01559"; E = "  -0.20%  " },
    @{ Row = 41; D = "1.005"; E = "  -1.05%  " },
    @{ Row = 42; D = "5.667"; E = "  +0.44%  " },
    @{ Row = 43; D = "0.8138"; E = "  -0.60%  " },
    @{ Row = 44; D = "99.27"; E = "  -0.54%  " },
    @{ Row = 45; D = "0.0₈122"; E = "  +6.04%  " },
    @{ Row = 46; D = "1.772.23"; E = "  -0.18%  " },
    @{ Row = 47; D = "0.4521"; E = "  -0.98%  " },
    @{ Row = 48; D = "1.009"; E = "  -0.68%  " },
    @{ Row = 49; D = "55.28"; E = "  +0.57%  " },
    @{ Row = 50; D = "0.05045"; E = "  +0.03%  " },
    @{ Row = 51; D = "1.007"; E = "  -0.77%  " }
)

# Preserve the existing (default) cell style while forcing these values
# to be stored as text, since several prices (e.g. "1.005") would
# otherwise be auto-coerced into numbers by Excel.
$origStyle = $ws.Range("B2").Style

foreach ($u in $updates) {
    $dCell = $ws.Range("D" + $u.Row)
    $dCell.NumberFormat = "@"
    $dCell.Value = $u.D
    $dCell.Style = $origStyle

    $eCell = $ws.Range("E" + $u.Row)
    $eCell.NumberFormat = "@"
    $eCell.Value = $u.E
    $eCell.Style = $origStyle
}
